# Lower case for name checking
# Append new artist/track name entries to the single-column "Names" list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @(
    "Shahin Loo",
    "Sasy",
    "Talk down",
    "Armin 2afm",
    "armin 2afm",
    "Armin zarei",
    "Armin Zarei",
    "Armin zarei"
)

# Find the first empty row below the existing data in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Range("A1").Value -eq $null) {
    $lastRow = 0
}

$row = $lastRow + 1
foreach ($name in $newNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

$ws.Range("A$row").Select()
